$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.1
$ws.Range("I2").Value = 3.75
$ws.Range("W2").Value = 6
$ws.Range("AF2").Value = 67
$ws.Range("AV2").Value = 67
$ws.Range("AY2").Value = 34

$ws.Rows("6").Delete()
